# "Generate Report for Archive"
# The localization report was regenerated: the in-flight status text moved on
# from "Ready for handoff" to "In Translation", and the now-shorter status
# values let the Status column narrow a bit in the refreshed report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# Re-fit the narrower status columns to their new (shorter) content.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
